$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text -like $needle) {
            return $para
        }
    }
    return $null
}

# --- Edit 1: Degree line (B.S. COMPUTER SCIENCE & MATH | SEATTLE UNIVERSITY) ---
$degreePara = Find-ParagraphByText $d "*B.S. COMPUTER SCIENCE*SEATTLE UNIVERSITY*"
if ($degreePara -eq $null) { throw "Degree paragraph not found" }
$pStart = $degreePara.Range.Start
$pEnd = $degreePara.Range.End
$degreeRange = $d.Range($pStart, $pEnd - 1)
$degreeXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">B.S. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">Major in </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:t>COMPUTER SCIENCE</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> &amp; Minor in PHILOSOPHY</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> | SEATTLE UNIVERSITY</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> ‘22</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:color w:val="191919"/><w:sz w:val="20"/></w:rPr><w:tab/><w:t xml:space="preserve"> 2020-CURRENT</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:color w:val="404040"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$degreeRange.InsertXML($degreeXml)

# --- Edit 2: GPA line (GPA: 3.33 -> GPA: 3.62, bold) ---
$gpaPara = Find-ParagraphByText $d "*GPA: 3.33*"
if ($gpaPara -eq $null) { throw "GPA paragraph not found" }
$gStart = $gpaPara.Range.Start
$gEnd = $gpaPara.Range.End
$gpaRange = $d.Range($gStart, $gEnd - 1)
$gpaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:color w:val="404040"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">GPA: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:bCs/><w:color w:val="404040"/><w:sz w:val="20"/></w:rPr><w:t>3.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:bCs/><w:color w:val="404040"/><w:sz w:val="20"/></w:rPr><w:t>62</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:color w:val="404040"/><w:sz w:val="20"/></w:rPr><w:tab/><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$gpaRange.InsertXML($gpaXml)

Write-Host "Degree paragraph now:" $degreePara.Range.Text
Write-Host "GPA paragraph now:" $gpaPara.Range.Text
